$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Existing rows 2-7: add column E submission labels ---
$ws.Range("E2").Value = "SubmissionLogBag"
$ws.Range("E3").Value = "SubmissionLogSimple"
$ws.Range("E4").Value = "SubmissionRFBag"
$ws.Range("E5").Value = "SubmissionRFSimple"
$ws.Range("E6").Value = "SubmissionCARTSimple"
$ws.Range("E7").Value = "SubmissionCARTBag"

# Row 6 (simple tree) previously had an empty B cell; it now gets a value
$ws.Range("B6").Value = 0.61426000000000003

# --- New block: rows 9-13 (row 8 left blank) ---
$ws.Range("A9").Value = "glm:WordCount+Weekday+Hour"
$ws.Range("B9").Value = 0.70491999999999999
$ws.Range("C9").Value = 0.83265310000000003
$ws.Range("D9").Value = 0.73701320000000003
$ws.Range("E9").Value = "SubmissionSimplestLog"

$ws.Range("A10").Value = "glm:WordCount+Weekday+Hour+NewsDesk"
$ws.Range("B10").Value = 0.88571
$ws.Range("C10").Value = 0.87704079999999995
$ws.Range("D10").Value = 0.89479390000000003
$ws.Range("E10").Value = "SubmissionSimplestLog"

$ws.Range("A11").Value = "glm:WordCount+Weekday+Hour+NewsDesk+SectionName"
$ws.Range("B11").Value = 0.89671000000000001
$ws.Range("C11").Value = 0.89744900000000005
$ws.Range("D11").Value = 0.91554760000000002
$ws.Range("E11").Value = "SubmissionSimplestLog"

$ws.Range("A12").Value = "RF: WordCount+NewsDesk+Hour+Weekday"
$ws.Range("B12").Value = 0.91107000000000005
$ws.Range("C12").Value = 0.89795919999999996
$ws.Range("D12").Value = 0.91921189999999997
$ws.Range("E12").Value = "SubmissionSimplestRF"

$ws.Range("A13").Value = "RF: WordCount+NewsDesk+Hour+Weekday+SectionName"
$ws.Range("C13").Value = 0.91326529999999995
$ws.Range("D13").Value = 0.93647460000000005
$ws.Range("E13").Value = "SubmissionSimplestRF"

# Custom number format (5 decimals) applied to the new "auc" column cells in rows 9-13
$ws.Range("C9:C13").NumberFormat = "0.00000_ "

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 52.0
$ws.Columns.Item(5).ColumnWidth = 22.0

# --- Selection matches the post-edit active cell in the source workbook ---
$ws.Range("D18").Select()
